$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 334
$ws.Range("F2").Value = 45200
$ws.Range("G2").Value = 30042
$ws.Range("H2").Value = 45231

# Row 3
$ws.Range("E3").Value = 29983
$ws.Range("F3").Value = 45170

# Row 4
$ws.Range("E4").Value = 29983
$ws.Range("F4").Value = 45170
$ws.Range("G4").Value = 30042
$ws.Range("H4").Value = 45231

# Row 5
$ws.Range("E5").Value = 30011
$ws.Range("F5").Value = 45200
$ws.Range("G5").Value = 30042
$ws.Range("H5").Value = 45231

# Row 6
$ws.Range("C6").Value = 442
$ws.Range("F6").Value = 45170
$ws.Range("G6").Value = 30042
$ws.Range("H6").Value = 45231

# Row 7
$ws.Range("E7").Value = 29952
$ws.Range("F7").Value = 45139
$ws.Range("G7").Value = 30042
$ws.Range("H7").Value = 45231

# Row 8
$ws.Range("D8").Value = 409
$ws.Range("E8").Value = 29983
$ws.Range("F8").Value = 45170
$ws.Range("H8").Value = 45231

# Row 9
$ws.Range("E9").Value = 29983
$ws.Range("F9").Value = 45170
$ws.Range("G9").Value = 30042
$ws.Range("H9").Value = 45231

# Row 10
$ws.Range("C10").Value = 500
$ws.Range("D10").Value = 482
$ws.Range("E10").Value = 29983
$ws.Range("F10").Value = 45170
$ws.Range("H10").Value = 45231

# Row 11
$ws.Range("E11").Value = 29983
$ws.Range("F11").Value = 45170
$ws.Range("G11").Value = 30042
$ws.Range("H11").Value = 45231

# Row 12
$ws.Range("C12").Value = 370
$ws.Range("D12").Value = 351
$ws.Range("F12").Value = 45170
$ws.Range("H12").Value = 45231

# Row 13
$ws.Range("C13").Value = 466
$ws.Range("F13").Value = 45170
$ws.Range("G13").Value = 30042
$ws.Range("H13").Value = 45231

# Row 14
$ws.Range("C14").Value = 423
$ws.Range("D14").Value = 395
$ws.Range("F14").Value = 45170
$ws.Range("H14").Value = 45231

# Row 15
$ws.Range("C15").Value = 383
$ws.Range("F15").Value = 45170
$ws.Range("G15").Value = 30011
$ws.Range("H15").Value = 45231

# Row 16
$ws.Range("C16").Value = 454
$ws.Range("D16").Value = 409
$ws.Range("F16").Value = 45170
$ws.Range("H16").Value = 45231

# Row 17
$ws.Range("C17").Value = 369
$ws.Range("D17").Value = 393
$ws.Range("F17").Value = 45170
$ws.Range("H17").Value = 45231

# Row 18
$ws.Range("E18").Value = 29983
$ws.Range("F18").Value = 45170
$ws.Range("G18").Value = 30042
$ws.Range("H18").Value = 45231

# Row 19
$ws.Range("D19").Value = 397
$ws.Range("E19").Value = 29983
$ws.Range("F19").Value = 45170
$ws.Range("H19").Value = 45231

# Row 20
$ws.Range("C20").Value = 480
$ws.Range("F20").Value = 45170
$ws.Range("G20").Value = 30042
$ws.Range("H20").Value = 45231

# Row 21
$ws.Range("E21").Value = 29983
$ws.Range("F21").Value = 45170
$ws.Range("G21").Value = 30042
$ws.Range("H21").Value = 45231

# Row 22
$ws.Range("D22").Value = 379
$ws.Range("E22").Value = 29983
$ws.Range("F22").Value = 45170
$ws.Range("H22").Value = 45231

# Row 23
$ws.Range("D23").Value = 216
$ws.Range("E23").Value = 29983
$ws.Range("F23").Value = 45170
$ws.Range("G23").Value = 35339
$ws.Range("H23").Value = 45231

# Row 24
$ws.Range("C24").Value = 393
$ws.Range("F24").Value = 45170
$ws.Range("G24").Value = 30042
$ws.Range("H24").Value = 45231

# Row 25
$ws.Range("E25").Value = 29983
$ws.Range("F25").Value = 45170
$ws.Range("G25").Value = 30042
$ws.Range("H25").Value = 45231

# Row 26
$ws.Range("C26").Value = 454
$ws.Range("D26").Value = 397
$ws.Range("F26").Value = 45200
$ws.Range("H26").Value = 45231

# Row 27
$ws.Range("C27").Value = 500
$ws.Range("D27").Value = 409
$ws.Range("E27").Value = 30011
$ws.Range("F27").Value = 45200
$ws.Range("H27").Value = 45231

# Row 28
$ws.Range("C28").Value = 358
$ws.Range("D28").Value = 366
$ws.Range("F28").Value = 45170
$ws.Range("H28").Value = 45231
